# Update correction and processed data from update_data
$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $b, $c, $d, $e, $f, $g) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
}

# --- Age Group sheet ---
$ws = $wb.Worksheets.Item("Age Group")
Set-Row $ws 2  50715 4353 2    8.390000000000001 7.94               0.08
Set-Row $ws 3  93646 9394 8    15.49              17.14              0.31
Set-Row $ws 4  93752 8689 18   15.51              15.85              0.6899999999999999
Set-Row $ws 5  90708 8993 46   15                 16.41              1.76
Set-Row $ws 6  97923 8421 133  16.2               15.36              5.1
Set-Row $ws 7  85840 6456 416  14.2               11.78              15.94
Set-Row $ws 8  53887 4069 643  8.91               7.42               24.64
Set-Row $ws 9  34878 4404 1344 5.77               8.029999999999999  51.49
Set-Row $ws 10 3286  34   0    0.54               0.06               0

# --- Gender sheet ---
$ws = $wb.Worksheets.Item("Gender")
Set-Row $ws 2 346093 28633 1270 57.24 52.24 48.66
Set-Row $ws 3 246653 25231 1286 40.79 46.03 49.27
Set-Row $ws 4 11889  949   54   1.97  1.73  2.07

# --- Race sheet ---
$ws = $wb.Worksheets.Item("Race")
Set-Row $ws 2 5064   753   13   0.84  1.37  0.5
Set-Row $ws 3 52411  6426  373  8.67  11.72 14.29
Set-Row $ws 4 81261  9425  347  13.44 17.19 13.3
Set-Row $ws 5 113387 12867 186  18.75 23.47 7.13
Set-Row $ws 6 352512 25342 1691 58.3  46.23 64.79000000000001

# --- Ethnicity sheet ---
$ws = $wb.Worksheets.Item("Ethnicity")
Set-Row $ws 2 20163  6102  51   3.33  11.13 1.95
Set-Row $ws 3 233526 18451 1373 38.62 33.66 52.61
Set-Row $ws 4 350946 30260 1186 58.04 55.21 45.44
